$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '27.550.37'
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '1.752.45'
$ws.Range('E3').NumberFormat = '@'
$ws.Range('E3').Value = '  -3.45%  '
$ws.Range('E4').NumberFormat = '@'
$ws.Range('E4').Value = '  +0.08%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '324.39'
$ws.Range('E5').NumberFormat = '@'
$ws.Range('E5').Value = '  -0.42%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '1.001'
$ws.Range('E6').NumberFormat = '@'
$ws.Range('E6').Value = '  +0.08%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.4478'
$ws.Range('E7').NumberFormat = '@'
$ws.Range('E7').Value = '  +2.52%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.3620'
$ws.Range('E8').NumberFormat = '@'
$ws.Range('E8').Value = '  -1.19%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.07511'
$ws.Range('E9').NumberFormat = '@'
$ws.Range('E9').Value = '  -2.06%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '42.14'
$ws.Range('E10').NumberFormat = '@'
$ws.Range('E10').Value = '  -5.84%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '1.105'
$ws.Range('E11').NumberFormat = '@'
$ws.Range('E11').Value = '  -3.11%  '
$ws.Range('E12').NumberFormat = '@'
$ws.Range('E12').Value = '  +0.10%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '20.68'
$ws.Range('E13').NumberFormat = '@'
$ws.Range('E13').Value = '  -5.83%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '6.049'
$ws.Range('E14').NumberFormat = '@'
$ws.Range('E14').Value = '  -4.02%  '
$ws.Range('E15').NumberFormat = '@'
$ws.Range('E15').Value = '  -4.07%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '1.754.30'
$ws.Range('E16').NumberFormat = '@'
$ws.Range('E16').Value = '  -3.64%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '92.85'
$ws.Range('E17').NumberFormat = '@'
$ws.Range('E17').Value = '  -2.21%  '
$ws.Range('E18').NumberFormat = '@'
$ws.Range('E18').Value = '  -1.36%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '0.06428'
$ws.Range('E19').NumberFormat = '@'
$ws.Range('E19').Value = '  -0.91%  '
$ws.Range('E20').NumberFormat = '@'
$ws.Range('E20').Value = '  +0.10%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '17.05'
$ws.Range('E21').NumberFormat = '@'
$ws.Range('E21').Value = '  -1.82%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '5.840'
$ws.Range('E22').NumberFormat = '@'
$ws.Range('E22').Value = '  -6.37%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '27.591.64'
$ws.Range('E23').NumberFormat = '@'
$ws.Range('E23').Value = '  -2.49%  '
$ws.Range('E24').NumberFormat = '@'
$ws.Range('E24').Value = '  -2.62%  '
$ws.Range('E25').NumberFormat = '@'
$ws.Range('E25').Value = '  -0.70%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '162.30'
$ws.Range('E26').NumberFormat = '@'
$ws.Range('E26').Value = '  +0.57%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '20.41'
$ws.Range('E27').NumberFormat = '@'
$ws.Range('E27').Value = '  -1.46%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '1.953.20'
$ws.Range('E28').NumberFormat = '@'
$ws.Range('E28').Value = '  -3.72%  '
$ws.Range('E29').NumberFormat = '@'
$ws.Range('E29').Value = '  -6.16%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '125.47'
$ws.Range('E30').NumberFormat = '@'
$ws.Range('E30').Value = '  -2.89%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '1.082'
$ws.Range('E31').NumberFormat = '@'
$ws.Range('E31').Value = '  -10.51%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '0.09021'
$ws.Range('E32').NumberFormat = '@'
$ws.Range('E32').Value = '  -1.26%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '3.639'
$ws.Range('E33').NumberFormat = '@'
$ws.Range('E33').Value = '  +2.26%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '5.546'
$ws.Range('E34').NumberFormat = '@'
$ws.Range('E34').Value = '  -7.69%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '12.02'
$ws.Range('E35').NumberFormat = '@'
$ws.Range('E35').Value = '  -7.30%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.02323'
$ws.Range('E36').NumberFormat = '@'
$ws.Range('E36').Value = '  -1.72%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '0.2096'
$ws.Range('E37').NumberFormat = '@'
$ws.Range('E37').Value = '  -3.54%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.6366'
$ws.Range('E38').NumberFormat = '@'
$ws.Range('E38').Value = '  -3.37%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.05974'
$ws.Range('E39').NumberFormat = '@'
$ws.Range('E39').Value = '  -3.82%  '
$ws.Range('B40').NumberFormat = '@'
$ws.Range('B40').Value = 'TrustWalletToken'
$ws.Range('C40').NumberFormat = '@'
$ws.Range('C40').Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '1.217'
$ws.Range('E40').NumberFormat = '@'
$ws.Range('E40').Value = '  +1.95%  '
$ws.Range('B41').NumberFormat = '@'
$ws.Range('B41').Value = 'InternetComputer(DFINITY)'
$ws.Range('C41').NumberFormat = '@'
$ws.Range('C41').Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '4.956'
$ws.Range('E41').NumberFormat = '@'
$ws.Range('E41').Value = '  -5.23%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '1.0000'
$ws.Range('E42').NumberFormat = '@'
$ws.Range('E42').Value = '  +0.12%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '1.386'
$ws.Range('E43').NumberFormat = '@'
$ws.Range('E43').Value = '  -2.88%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '7.807'
$ws.Range('E44').NumberFormat = '@'
$ws.Range('E44').Value = '  -3.18%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '13.23'
$ws.Range('E45').NumberFormat = '@'
$ws.Range('E45').Value = '  -4.22%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '0.5884'
$ws.Range('E46').NumberFormat = '@'
$ws.Range('E46').Value = '  -3.55%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '3.712'
$ws.Range('E47').NumberFormat = '@'
$ws.Range('E47').Value = '  -0.62%  '
$ws.Range('E48').NumberFormat = '@'
$ws.Range('E48').Value = '  -2.78%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '121.26'
$ws.Range('E49').NumberFormat = '@'
$ws.Range('E49').Value = '  -3.33%  '
$ws.Range('E50').NumberFormat = '@'
$ws.Range('E50').Value = '  +0.30%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.06865'
$ws.Range('E51').NumberFormat = '@'
$ws.Range('E51').Value = '  -1.84%  '
